$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1): columns C, D, E
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Update data rows 2 through 25: column C becomes the "f__UBA660" text (same as D),
# column D stays "f__UBA660", and column E becomes numeric 1
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 3).Value = "f__UBA660"
    $ws.Cells.Item($r, 4).Value = "f__UBA660"
    $ws.Cells.Item($r, 5).Value = 1
}
